$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 2).Value = 1.732704007046913
$ws.Cells.Item(1, 3).Value = 0.1945306715051764
$ws.Cells.Item(1, 4).Value = -0.672838158291254
$ws.Cells.Item(1, 5).Value = 0.7034277224914169
$ws.Cells.Item(1, 6).Value = 1.570796292848413

$ws.Cells.Item(2, 2).Value = 1.734449897926843
$ws.Cells.Item(2, 3).Value = 0.1944889771920187
$ws.Cells.Item(2, 4).Value = -0.6729541981860235
$ws.Cells.Item(2, 5).Value = 0.7033533778325018
$ws.Cells.Item(2, 6).Value = 1.570796292708984

$ws.Cells.Item(3, 2).Value = 1.745338072090327
$ws.Cells.Item(3, 3).Value = 0.1942289523826182
$ws.Cells.Item(3, 4).Value = -0.6736778760508239
$ws.Cells.Item(3, 5).Value = 0.7028897305314603
$ws.Cells.Item(3, 6).Value = 1.570796291839441

$ws.Cells.Item(4, 2).Value = 1.771088086074367
$ws.Cells.Item(4, 3).Value = 0.1936140060768625
$ws.Cells.Item(4, 4).Value = -0.6753893397782359
$ws.Cells.Item(4, 5).Value = 0.7017932267185878
$ws.Cells.Item(4, 6).Value = 1.570796289783013

$ws.Cells.Item(5, 2).Value = 1.814176913182388
$ws.Cells.Item(5, 3).Value = 0.1925849846920256
$ws.Cells.Item(5, 4).Value = -0.6782532203453618
$ws.Cells.Item(5, 5).Value = 0.6999583903085791
$ws.Cells.Item(5, 6).Value = 1.570796286341886

$ws.Cells.Item(6, 2).Value = 1.87437937402317
$ws.Cells.Item(6, 3).Value = 0.1911472658265371
$ws.Cells.Item(6, 4).Value = -0.6822545512997377
$ws.Cells.Item(6, 5).Value = 0.6973948100364618
$ws.Cells.Item(6, 6).Value = 1.570796281534042

$ws.Cells.Item(7, 2).Value = 1.949308567049777
$ws.Cells.Item(7, 3).Value = 0.1893578520237512
$ws.Cells.Item(7, 4).Value = -0.6872346882452469
$ws.Cells.Item(7, 5).Value = 0.6942041264935301
$ws.Cells.Item(7, 6).Value = 1.570796275550103

$ws.Cells.Item(8, 2).Value = 2.03495629909848
$ws.Cells.Item(8, 3).Value = 0.1873124645357151
$ws.Cells.Item(8, 4).Value = -0.6929272283280338
$ws.Cells.Item(8, 5).Value = 0.6905570191632778
$ws.Cells.Item(8, 6).Value = 1.570796268710169

$ws.Cells.Item(9, 2).Value = 2.126233515927697
$ws.Cells.Item(9, 3).Value = 0.1851326370869383
$ws.Cells.Item(9, 4).Value = -0.6989939297224161
$ws.Cells.Item(9, 5).Value = 0.686670193457332
$ws.Cells.Item(9, 6).Value = 1.570796261420657

$ws.Cells.Item(10, 2).Value = 2.217510732756915
$ws.Cells.Item(10, 3).Value = 0.1829528096381616
$ws.Cells.Item(10, 4).Value = -0.7050606311167984
$ws.Cells.Item(10, 5).Value = 0.6827833677513863
$ws.Cells.Item(10, 6).Value = 1.570796254131145

$ws.Cells.Item(11, 2).Value = 2.303158464805618
$ws.Cells.Item(11, 3).Value = 0.1809074221501255
$ws.Cells.Item(11, 4).Value = -0.7107531711995853
$ws.Cells.Item(11, 5).Value = 0.679136260421134
$ws.Cells.Item(11, 6).Value = 1.57079624729121

$ws.Cells.Item(12, 2).Value = 2.378087657832224
$ws.Cells.Item(12, 3).Value = 0.1791180083473395
$ws.Cells.Item(12, 4).Value = -0.7157333081450945
$ws.Cells.Item(12, 5).Value = 0.6759455768782023
$ws.Cells.Item(12, 6).Value = 1.570796241307272

$ws.Cells.Item(13, 2).Value = 2.438290118673007
$ws.Cells.Item(13, 3).Value = 0.1776802894818511
$ws.Cells.Item(13, 4).Value = -0.7197346390994703
$ws.Cells.Item(13, 5).Value = 0.673381996606085
$ws.Cells.Item(13, 6).Value = 1.570796236499428

$ws.Cells.Item(14, 2).Value = 2.481378945781028
$ws.Cells.Item(14, 3).Value = 0.1766512680970142
$ws.Cells.Item(14, 4).Value = -0.7225985196665963
$ws.Cells.Item(14, 5).Value = 0.6715471601960763
$ws.Cells.Item(14, 6).Value = 1.570796233058301

$ws.Cells.Item(15, 2).Value = 2.507128959765068
$ws.Cells.Item(15, 3).Value = 0.1760363217912584
$ws.Cells.Item(15, 4).Value = -0.7243099833940083
$ws.Cells.Item(15, 5).Value = 0.6704506563832038
$ws.Cells.Item(15, 6).Value = 1.570796231001873

$ws.Cells.Item(16, 2).Value = 2.518017133928552
$ws.Cells.Item(16, 3).Value = 0.1757762969818579
$ws.Cells.Item(16, 4).Value = -0.7250336612588086
$ws.Cells.Item(16, 5).Value = 0.6699870090821622
$ws.Cells.Item(16, 6).Value = 1.57079623013233

$ws.Cells.Item(17, 2).Value = 2.519763024808483
$ws.Cells.Item(17, 3).Value = 0.1757346026687003
$ws.Cells.Item(17, 4).Value = -0.7251497011535781
$ws.Cells.Item(17, 5).Value = 0.6699126644232472
$ws.Cells.Item(17, 6).Value = 1.570796229992901
